# Update the return_status / note pair for row 3 (sn202508 / 321903):
#   return_status: NCC_DA_NHAN_TRA_HANG -> TIEP_NHAN_TRA_HANG
#   note:          Viettien đã đồng ý hoàn hàng -> Hương xác nhận đồng ý nhận lại hàng
# Update the return_status / note pair for row 4 (sn202512 / 321895):
#   return_status: NCC_DA_DOI_HANG -> DA_TRA_HANG
#   note:          Katinat đã đổi lý mới cho anh H -> Hương đã nhận trả hàng thành công

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "TIEP_NHAN_TRA_HANG"
$ws.Range("D3").Value = "Hương xác nhận đồng ý nhận lại hàng"

$ws.Range("C4").Value = "DA_TRA_HANG"
$ws.Range("D4").Value = "Hương đã nhận trả hàng thành công"

# Move the active selection to where the author left it when saving.
$ws.Range("D8").Select()
